# "Generate Report for Archive"
#
# The localization status changes from "Ready for handoff" to
# "In Translation" for the tracked file, on all three sheets that
# reference it ("Overview" E2/F2, "zh-cn" C2, "de-de" C2). Because
# "Ready for handoff" / "In Translation" is stored once as a shared
# string, updating the cell Value on every occurrence collapses back
# onto a single shared-string entry, exactly like the source diff.
#
# The status columns that held the long status text are narrower once
# the shorter "In Translation" text is in place (mirrors an
# autofit-after-edit in the original workbook), so their column widths
# shrink correspondingly.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

# Narrow the status columns to match the new (shorter) content width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
